$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SEER Survival")

# --- New Weibull calculation table (rows 33-36) ---
# Shared-string insertion order matters for byte-identical output, so cells
# are written in the same order the original author entered them:
# title -> S(5)/S(10) headers -> Group label -> p5/p10/shape/scale headers -> url
$ws.Range("B33").Value = "Weibull calculation"
$ws.Range("B34").Value = "S(5)"
$ws.Range("C34").Value = "S(10)"
$ws.Range("A33").Value = "Group"
$ws.Range("D34").Value = "p5"
$ws.Range("E34").Value = "p10"
$ws.Range("F34").Value = "shape"
$ws.Range("G34").Value = "scale"
$ws.Range("D33").Value = "http://www.johndcook.com/quantiles_parameters.pdf"

$ws.Range("A35").Value = "Reg+Dist"
$ws.Range("B35").Formula = "=D16"
$ws.Range("C35").Formula = "=D17"
$ws.Range("D35").Formula = "=1-B35"
$ws.Range("E35").Formula = "=1-C35"
$ws.Range("F35").Formula = "=(LN(-LN(1-E35))-LN(-LN(1-D35)))/(LN(10)-LN(5))"
$ws.Range("G35").Formula = "=5/(-LN(1-D35))^(1/F35)"

$ws.Range("A36").Value = "Local"
$ws.Range("B36").Formula = "=D4"
$ws.Range("C36").Formula = "=D5"
$ws.Range("D36").Formula = "=1-B36"
$ws.Range("E36").Formula = "=1-C36"
$ws.Range("F36").Formula = "=(LN(-LN(1-E36))-LN(-LN(1-D36)))/(LN(10)-LN(5))"
$ws.Range("G36").Formula = "=5/(-LN(1-D36))^(1/F36)"

# --- Sheet activation / selection (SEER Survival becomes the active tab) ---
$ws.Activate()
$ws.Range("D33").Select()

Write-Host "Done"
